# Insert a new weekly record at the top of the data table (row 8), which
# pushes the existing rows 8-23 down to rows 9-24 (old row 23 ends up
# duplicated as the new row 24, matching a rolling weekly window).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(8).Insert()

# Populate the new row 8 with this week's record.
$ws.Cells.Item(8, 1).Value  = 11
$ws.Cells.Item(8, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(8, 3).Value  = "Bíobío"
$ws.Cells.Item(8, 4).Value  = 44447
$ws.Cells.Item(8, 5).Value  = 8
$ws.Cells.Item(8, 6).Value  = 100112013
$ws.Cells.Item(8, 7).Value  = "Alcachofa"
$ws.Cells.Item(8, 8).Value  = "Española"
$ws.Cells.Item(8, 9).Value  = "Primera"
$ws.Cells.Item(8, 10).Value = 100
$ws.Cells.Item(8, 11).Value = 14000
$ws.Cells.Item(8, 12).Value = 15000
$ws.Cells.Item(8, 13).Value = 14500
$ws.Cells.Item(8, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(8, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(8, 16).Value = 483
$ws.Cells.Item(8, 17).Value = 30
$ws.Cells.Item(8, 18).Value = "Hortaliza"
